$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = "Opinioni positive sul biglietto di gratta e vinci"

# Update row 3
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Opinioni negative sul biglietto della lotteria"

# Delete rows 4 through 19 (the remaining data rows)
$ws.Range("A4:C19").EntireRow.Delete()
